$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# but after deleting row 26 it becomes row 27).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# After the two row deletions, the remaining rows have shifted up by two.
# Update the "B" column (column C) values that changed for the surviving
# rows (SC 101 gains a value, SC 119 and SC 193 lose theirs).
$ws.Range("C27").Value = 10
$ws.Range("C29").ClearContents()
$ws.Range("C32").ClearContents()
